$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 44
$ws1.Range("F13").Value = 1831
$ws1.Range("F18").Value = 6277
$ws1.Range("F19").Value = 242
$ws1.Range("F23").Value = 883
$ws1.Range("F27").Value = 2473
$ws1.Range("F33").Value = 1318
$ws1.Range("F37").Value = 32
$ws1.Range("F39").Value = 1508
$ws1.Range("F41").Value = 1465
$ws1.Range("F42").Value = 94

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 219
$ws2.Range("F11").Value = 168
$ws2.Range("F17").Value = 340

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 926
$ws3.Range("F4").Value = 262
$ws3.Range("F6").Value = 46

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 926
$ws4.Range("F7").Value = 262
$ws4.Range("F10").Value = 46
$ws4.Range("F18").Value = 44
$ws4.Range("F23").Value = 168
$ws4.Range("F27").Value = 6277
$ws4.Range("F28").Value = 242
$ws4.Range("F34").Value = 2474
$ws4.Range("F37").Value = 1318
$ws4.Range("F40").Value = 340
$ws4.Range("F45").Value = 32
$ws4.Range("F48").Value = 1508
$ws4.Range("F51").Value = 94
